# ------------------------------------------------------------------
# Fixed update to excel issue
#   - Rename "Requested quantity" header on the two existing sheets
#   - Add a new "PO Forecast" sheet with forecast data
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" headers
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# ------------------------------------------------------------------
# Add the new "PO Forecast" sheet after the last existing sheet
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Match page margins used by the rest of the workbook
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

# Headers
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(45494.99999999999, 31, 22.03718455769523, 39.70079724793253),
    @(45508.99999999999, 32, 23.39952060742251, 40.67797419175051),
    @(45557.99999999999, 37, 27.95081372525929, 46.00511254890422),
    @(45592.99999999999, 40, 31.57140507106969, 47.99497272522432),
    @(45599.99999999999, 40, 31.7260109241368,  48.79486829196532),
    @(45606.99999999999, 41, 32.54834921455443, 49.5573920307493),
    @(45613.99999999999, 42, 32.69082133356536, 50.45203868549078),
    @(45620.99999999999, 42, 32.95523561734601, 50.58569620810874),
    @(45627.99999999999, 43, 33.81782622447736, 51.72396146141195),
    @(45634.99999999999, 44, 34.78405188211401, 52.29838503459268),
    @(45641.99999999999, 44, 36.10708869442159, 52.74422058119757),
    @(45648.99999999999, 45, 35.89958139399354, 53.75212096037978),
    @(45655.99999999999, 46, 36.55085717698107, 54.55824373889148)
)

$row = 2
foreach ($r in $data) {
    $ws3.Cells.Item($row, 1).Value = $r[0]
    $ws3.Cells.Item($row, 2).Value = $r[1]
    $ws3.Cells.Item($row, 3).Value = $r[2]
    $ws3.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# ------------------------------------------------------------------
# Re-use the existing style definitions (bold/centered header with
# border, and the YYYY-MM-DD HH:MM:SS date format) instead of minting
# new ones, by copy/pasting formats from the already-styled cells on
# the "Weekly Quantity" sheet.
# ------------------------------------------------------------------
$ws1.Range("B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws3.Range("A2:A14").PasteSpecial(-4122)

$ws3.Range("A1").Select() | Out-Null
$excel.CutCopyMode = 0
